# Swap the "Periodo Mora" (E) and "Valor Mora" (F) values between
# rows 16 and 17 on the only worksheet ("Hoja1").
#
# Before:  E16=2201  F16=36341      After:  E16=2112  F16=18170
#          E17=2112  F17=18170              E17=2201  F17=36341

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("E16").Value = "2112"
$ws.Range("F16").Value = 18170

$ws.Range("E17").Value = "2201"
$ws.Range("F17").Value = 36341
